$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efnb2"
$ws.Range("C2").Value = "Ephb1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 36.899643
$ws.Range("H2").Value = 110.698929
$ws.Range("I2").Value = 0.7238945645409351
$ws.Range("J2").Value = 0.7238945645409351
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.471191666666667
$ws.Range("N2").Value = 4.413575
$ws.Range("O2").Value = 0.6447353255635294
$ws.Range("P2").Value = 0.6447353255635294
$ws.Range("Q2").Value = 54.286447284575
$ws.Range("R2").Value = 488.578025561175
$ws.Range("S2").Value = 0.4667203977429692
$ws.Range("T2").Value = 0.4667203977429692

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efnb2"
$ws.Range("C3").Value = "Ephb1"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 36.899643
$ws.Range("H3").Value = 110.698929
$ws.Range("I3").Value = 0.7238945645409351
$ws.Range("J3").Value = 0.7238945645409351
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.8106620000000001
$ws.Range("N3").Value = 2.431986
$ws.Range("O3").Value = 0.3552646744364706
$ws.Range("P3").Value = 0.3552646744364706
$ws.Range("Q3").Value = 29.913138393666
$ws.Range("R3").Value = 269.218245542994
$ws.Range("S3").Value = 0.257174166797966
$ws.Range("T3").Value = 0.257174166797966

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Efnb2"
$ws.Range("C4").Value = "Ephb1"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.374819
$ws.Range("H4").Value = 10.124457
$ws.Range("I4").Value = 0.0662069584361419
$ws.Range("J4").Value = 0.0662069584361419
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 1.471191666666667
$ws.Range("N4").Value = 4.413575
$ws.Range("O4").Value = 0.6447353255635294
$ws.Range("P4").Value = 0.6447353255635294
$ws.Range("Q4").Value = 4.965005589308333
$ws.Range("R4").Value = 44.68505030377499
$ws.Range("S4").Value = 0.04268596490189701
$ws.Range("T4").Value = 0.04268596490189701

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efnb2"
$ws.Range("C5").Value = "Ephb1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.374819
$ws.Range("H5").Value = 10.124457
$ws.Range("I5").Value = 0.0662069584361419
$ws.Range("J5").Value = 0.0662069584361419
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.8106620000000001
$ws.Range("N5").Value = 2.431986
$ws.Range("O5").Value = 0.3552646744364706
$ws.Range("P5").Value = 0.3552646744364706
$ws.Range("Q5").Value = 2.735837520178
$ws.Range("R5").Value = 24.622537681602
$ws.Range("S5").Value = 0.0235209935342449
$ws.Range("T5").Value = 0.0235209935342449

$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Efnb2"
$ws.Range("C6").Value = "Ephb1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 10.699319
$ws.Range("H6").Value = 32.097957
$ws.Range("I6").Value = 0.2098984770229228
$ws.Range("J6").Value = 0.2098984770229228
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 1.471191666666667
$ws.Range("N6").Value = 4.413575
$ws.Range("O6").Value = 0.6447353255635294
$ws.Range("P6").Value = 0.6447353255635294
$ws.Range("Q6").Value = 15.74074895180834
$ws.Range("R6").Value = 141.666740566275
$ws.Range("S6").Value = 0.1353289629186631
$ws.Range("T6").Value = 0.1353289629186631

$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Efnb2"
$ws.Range("C7").Value = "Ephb1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 10.699319
$ws.Range("H7").Value = 32.097957
$ws.Range("I7").Value = 0.2098984770229228
$ws.Range("J7").Value = 0.2098984770229228
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.8106620000000001
$ws.Range("N7").Value = 2.431986
$ws.Range("O7").Value = 0.3552646744364706
$ws.Range("P7").Value = 0.3552646744364706
$ws.Range("Q7").Value = 8.673531339178002
$ws.Range("R7").Value = 78.06178205260201
$ws.Range("S7").Value = 0.07456951410425969
$ws.Range("T7").Value = 0.07456951410425969
